# PowerShell COM-interop script to update cryptos.xlsx price/volume data
# Generated to match the target commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap rows 37 and 38: PEPE moves to row 37, OKB moves to row 38, with updated values ---
$ws.Range("B37").Value = "PEPE"
$ws.Range("C37").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0₃0761"
$ws.Range("E37").Value = "  -1.42%  "

$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "52.76"
$ws.Range("E38").Value = "  -4.89%  "

# --- Update Price (D) and Volume(1h) (E) values for all other rows ---
$ws.Range("D2").Value = "63.429.63"
$ws.Range("E2").Value = "  -4.05%  "
$ws.Range("D3").Value = "3.123.78"
$ws.Range("E3").Value = "  -3.90%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "607.04"
$ws.Range("E5").Value = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.60"
$ws.Range("E6").Value = "  -7.89%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("D8").Value = "3.118.02"
$ws.Range("E8").Value = "  -4.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("E9").Value = "  -3.47%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.150"
$ws.Range("E10").Value = "  -6.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.25"
$ws.Range("E11").Value = "  -7.31%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.468"
$ws.Range("E12").Value = "  -4.36%  "
$ws.Range("E13").Value = "  -5.87%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "35.16"
$ws.Range("E14").Value = "  -8.21%  "
$ws.Range("D15").Value = "3.627.75"
$ws.Range("E15").Value = "  -4.01%  "
$ws.Range("E16").Value = "  +1.87%  "
$ws.Range("D17").Value = "63.546.05"
$ws.Range("E17").Value = "  -3.85%  "
$ws.Range("D18").Value = "3.118.76"
$ws.Range("E18").Value = "  -4.07%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.79"
$ws.Range("E19").Value = "  -5.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "474.23"
$ws.Range("E20").Value = "  -4.32%  "
$ws.Range("E21").Value = "  -4.62%  "
$ws.Range("E22").Value = "  -5.29%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "7.75"
$ws.Range("E23").Value = "  -2.55%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.42"
$ws.Range("E24").Value = "  -7.05%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.54"
$ws.Range("E25").Value = "  -3.52%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.78"
$ws.Range("E27").Value = "  -7.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.38"
$ws.Range("E28").Value = "  -6.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.118"
$ws.Range("E29").Value = "  -10.72%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.82"
$ws.Range("E30").Value = "  -1.98%  "
$ws.Range("E31").Value = "  -11.90%  "
$ws.Range("E32").Value = "  +0.05%  "
$ws.Range("E33").Value = "  -5.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "26.13"
$ws.Range("E34").Value = "  -5.55%  "
$ws.Range("E35").Value = "  -1.64%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.91"
$ws.Range("E36").Value = "  -7.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "452.28"
$ws.Range("E39").Value = "  -7.48%  "
$ws.Range("E40").Value = "  -13.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0391"
$ws.Range("E41").Value = "  -6.40%  "
$ws.Range("E42").Value = "  -9.17%  "
$ws.Range("E43").Value = "  -4.49%  "
$ws.Range("D44").Value = "2.837.08"
$ws.Range("E44").Value = "  -5.05%  "
$ws.Range("E45").Value = "  -10.26%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.263"
$ws.Range("E46").Value = "  -8.74%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.40"
$ws.Range("E47").Value = "  -1.84%  "
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "26.02"
$ws.Range("E49").Value = "  -8.00%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.113"
$ws.Range("E50").Value = "  -4.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "118.88"
$ws.Range("E51").Value = "  -2.39%  "

Write-Host "Crypto data updated successfully"
